$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q3" right before the existing "2022-Q2"
#    sheet. This shifts 2022-Q2 .. 2020-Q4 one position to the right, just
#    like the target workbook (sheet2 .. sheet8 become sheet3 .. sheet9).
# ---------------------------------------------------------------------------
$existingQ2 = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($existingQ2)
$q3.Name = "2022-Q3"

$summary = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------------
# 2. 总计 (summary) sheet: insert a new row 2 for the 2022-Q3 totals,
#    pushing the existing quarters down by one row (2022-Q2 -> row3, ...,
#    2020-Q4 -> row9).
# ---------------------------------------------------------------------------
$summary.Rows.Item(2).Insert()
$summary.Range("A2:D2").ClearFormats()

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 45
$summary.Cells.Item(2, 4).Value = 30.31

# Re-apply the index-column header style (bold/border/center, same as the
# untouched A3:A9 cells) to the new A2 cell.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. Populate the new "2022-Q3" sheet with the fund holdings table.
# ---------------------------------------------------------------------------
$headerSrc = $summary.Range("B1:D1")
$headerSrc.Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

$q3.Cells.Item(1, 2).Value = "基金代码"
$q3.Cells.Item(1, 3).Value = "基金名称"
$q3.Cells.Item(1, 4).Value = "基金规模"
$q3.Cells.Item(1, 5).Value = "股票总仓位"
$q3.Cells.Item(1, 6).Value = "仓位占比"
$q3.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q3.Cells.Item(1, 8).Value = "仓位排名"

# Force columns B (fund code) and D:G (text-formatted numbers) to be
# stored as text, matching the source data (t="inlineStr").
$q3.Range("B2:B46").NumberFormat = "@"
$q3.Range("D2:G46").NumberFormat = "@"

$arr = New-Object 'object[,]' 45,8
$arr[0,0] = 0
$arr[0,1] = "512880"
$arr[0,2] = "国泰中证全指证券公司ETF"
$arr[0,3] = "289.27"
$arr[0,4] = "99.95"
$arr[0,5] = "3.47"
$arr[0,6] = "10.0377"
$arr[0,7] = 6
$arr[1,0] = 1
$arr[1,1] = "512000"
$arr[1,2] = "华宝中证全指证券公司ETF"
$arr[1,3] = "215.91"
$arr[1,4] = "99.89"
$arr[1,5] = "3.50"
$arr[1,6] = "7.5568"
$arr[1,7] = 6
$arr[2,0] = 2
$arr[2,1] = "512900"
$arr[2,2] = "南方中证全指证券公司ETF"
$arr[2,3] = "78.74"
$arr[2,4] = "99.98"
$arr[2,5] = "3.51"
$arr[2,6] = "2.7638"
$arr[2,7] = 6
$arr[3,0] = 3
$arr[3,1] = "159841"
$arr[3,2] = "天弘中证全指证券公司ETF"
$arr[3,3] = "45.70"
$arr[3,4] = "99.94"
$arr[3,5] = "3.51"
$arr[3,6] = "1.6041"
$arr[3,7] = 6
$arr[4,0] = 4
$arr[4,1] = "512070"
$arr[4,2] = "易方达沪深300非银行金融ETF"
$arr[4,3] = "38.64"
$arr[4,4] = "99.35"
$arr[4,5] = "2.87"
$arr[4,6] = "1.1090"
$arr[4,7] = 9
$arr[5,0] = 5
$arr[5,1] = "159993"
$arr[5,2] = "鹏华国证证券龙头ETF"
$arr[5,3] = "13.39"
$arr[5,4] = "98.06"
$arr[5,5] = "7.11"
$arr[5,6] = "0.9520"
$arr[5,7] = 5
$arr[6,0] = 6
$arr[6,1] = "161720"
$arr[6,2] = "招商中证全指证券公司指数（LOF）A"
$arr[6,3] = "22.28"
$arr[6,4] = "94.50"
$arr[6,5] = "3.30"
$arr[6,6] = "0.7352"
$arr[6,7] = 6
$arr[7,0] = 7
$arr[7,1] = "501016"
$arr[7,2] = "国泰中证申万证券行业指数（LOF）A"
$arr[7,3] = "18.10"
$arr[7,4] = "93.42"
$arr[7,5] = "3.31"
$arr[7,6] = "0.5991"
$arr[7,7] = 6
$arr[8,0] = 8
$arr[8,1] = "163113"
$arr[8,2] = "申万菱信中证申万证券行业指数（LOF）A"
$arr[8,3] = "16.20"
$arr[8,4] = "93.19"
$arr[8,5] = "3.30"
$arr[8,6] = "0.5346"
$arr[8,7] = 6
$arr[9,0] = 9
$arr[9,1] = "001553"
$arr[9,2] = "天弘中证证券保险指数型发起式 C"
$arr[9,3] = "17.93"
$arr[9,4] = "94.82"
$arr[9,5] = "2.50"
$arr[9,6] = "0.4482"
$arr[9,7] = 9
$arr[10,0] = 10
$arr[10,1] = "161027"
$arr[10,2] = "富国中证全指证券公司指数A"
$arr[10,3] = "12.47"
$arr[10,4] = "94.32"
$arr[10,5] = "3.30"
$arr[10,6] = "0.4115"
$arr[10,7] = 6
$arr[11,0] = 11
$arr[11,1] = "502010"
$arr[11,2] = "易方达证券公司指数（LOF）A"
$arr[11,3] = "12.08"
$arr[11,4] = "94.58"
$arr[11,5] = "3.31"
$arr[11,6] = "0.3998"
$arr[11,7] = 6
$arr[12,0] = 12
$arr[12,1] = "160633"
$arr[12,2] = "鹏华中证全指证券公司指数（LOF）A"
$arr[12,3] = "11.96"
$arr[12,4] = "94.06"
$arr[12,5] = "3.29"
$arr[12,6] = "0.3935"
$arr[12,7] = 6
$arr[13,0] = 13
$arr[13,1] = "501048"
$arr[13,2] = "汇添富中证全指证券公司指数（LOF）C"
$arr[13,3] = "9.31"
$arr[13,4] = "93.53"
$arr[13,5] = "3.29"
$arr[13,6] = "0.3063"
$arr[13,7] = 6
$arr[14,0] = 14
$arr[14,1] = "001552"
$arr[14,2] = "天弘中证证券保险指数型发起式 A"
$arr[14,3] = "12.15"
$arr[14,4] = "94.82"
$arr[14,5] = "2.50"
$arr[14,6] = "0.3038"
$arr[14,7] = 9
$arr[15,0] = 15
$arr[15,1] = "515010"
$arr[15,2] = "华夏中证全指证券公司ETF"
$arr[15,3] = "8.47"
$arr[15,4] = "99.66"
$arr[15,5] = "3.48"
$arr[15,6] = "0.2948"
$arr[15,7] = 6
$arr[16,0] = 16
$arr[16,1] = "160625"
$arr[16,2] = "鹏华中证800证券保险指数（LOF）A"
$arr[16,3] = "9.50"
$arr[16,4] = "94.43"
$arr[16,5] = "2.52"
$arr[16,6] = "0.2394"
$arr[16,7] = 9
$arr[17,0] = 17
$arr[17,1] = "160516"
$arr[17,2] = "博时中证全指证券公司指数"
$arr[17,3] = "7.00"
$arr[17,4] = "93.76"
$arr[17,5] = "3.28"
$arr[17,6] = "0.2296"
$arr[17,7] = 6
$arr[18,0] = 18
$arr[18,1] = "012044"
$arr[18,2] = "鹏华中证全指证券公司指数（LOF）C"
$arr[18,3] = "5.89"
$arr[18,4] = "94.06"
$arr[18,5] = "3.29"
$arr[18,6] = "0.1938"
$arr[18,7] = 6
$arr[19,0] = 19
$arr[19,1] = "501047"
$arr[19,2] = "汇添富中证全指证券公司指数（LOF）A"
$arr[19,3] = "5.58"
$arr[19,4] = "93.53"
$arr[19,5] = "3.29"
$arr[19,6] = "0.1836"
$arr[19,7] = 6
$arr[20,0] = 20
$arr[20,1] = "159842"
$arr[20,2] = "银华中证全指证券公司ETF"
$arr[20,3] = "5.01"
$arr[20,4] = "98.00"
$arr[20,5] = "3.43"
$arr[20,6] = "0.1718"
$arr[20,7] = 6
$arr[21,0] = 21
$arr[21,1] = "515560"
$arr[21,2] = "建信中证全指证券公司ETF"
$arr[21,3] = "3.98"
$arr[21,4] = "98.63"
$arr[21,5] = "3.54"
$arr[21,6] = "0.1409"
$arr[21,7] = 6
$arr[22,0] = 22
$arr[22,1] = "160419"
$arr[22,2] = "华安中证证券公司A"
$arr[22,3] = "3.88"
$arr[22,4] = "94.43"
$arr[22,5] = "3.29"
$arr[22,6] = "0.1277"
$arr[22,7] = 6
$arr[23,0] = 23
$arr[23,1] = "502053"
$arr[23,2] = "长盛中证全指证券公司指数（LOF）"
$arr[23,3] = "3.93"
$arr[23,4] = "93.15"
$arr[23,5] = "3.22"
$arr[23,6] = "0.1265"
$arr[23,7] = 6
$arr[24,0] = 24
$arr[24,1] = "012874"
$arr[24,2] = "易方达证券公司指数（LOF）C"
$arr[24,3] = "2.22"
$arr[24,4] = "94.58"
$arr[24,5] = "3.31"
$arr[24,6] = "0.0735"
$arr[24,7] = 6
$arr[25,0] = 25
$arr[25,1] = "512570"
$arr[25,2] = "易方达中证全指证券公司ETF"
$arr[25,3] = "1.92"
$arr[25,4] = "98.99"
$arr[25,5] = "3.47"
$arr[25,6] = "0.0666"
$arr[25,7] = 6
$arr[26,0] = 26
$arr[26,1] = "515850"
$arr[26,2] = "富国中证全指证券公司ETF"
$arr[26,3] = "1.67"
$arr[26,4] = "99.74"
$arr[26,5] = "3.52"
$arr[26,6] = "0.0588"
$arr[26,7] = 6
$arr[27,0] = 27
$arr[27,1] = "510200"
$arr[27,2] = "汇安上证证券ETF"
$arr[27,3] = "0.67"
$arr[27,4] = "95.06"
$arr[27,5] = "5.52"
$arr[27,6] = "0.0370"
$arr[27,7] = 5
$arr[28,0] = 28
$arr[28,1] = "159848"
$arr[28,2] = "国联安中证全指证券公司ETF"
$arr[28,3] = "0.96"
$arr[28,4] = "96.87"
$arr[28,5] = "3.38"
$arr[28,6] = "0.0324"
$arr[28,7] = 6
$arr[29,0] = 29
$arr[29,1] = "515630"
$arr[29,2] = "鹏华中证800证券保险ETF"
$arr[29,3] = "1.11"
$arr[29,4] = "95.27"
$arr[29,5] = "2.55"
$arr[29,6] = "0.0283"
$arr[29,7] = 9
$arr[30,0] = 30
$arr[30,1] = "516730"
$arr[30,2] = "浦银安盛中证证券公司30ETF"
$arr[30,3] = "0.60"
$arr[30,4] = "97.43"
$arr[30,5] = "4.28"
$arr[30,6] = "0.0257"
$arr[30,7] = 6
$arr[31,0] = 31
$arr[31,1] = "013276"
$arr[31,2] = "富国中证全指证券公司指数C"
$arr[31,3] = "0.57"
$arr[31,4] = "94.32"
$arr[31,5] = "3.30"
$arr[31,6] = "0.0188"
$arr[31,7] = 6
$arr[32,0] = 32
$arr[32,1] = "012606"
$arr[32,2] = "西藏东财中证证券保险领先指数C"
$arr[32,3] = "0.61"
$arr[32,4] = "94.99"
$arr[32,5] = "2.98"
$arr[32,6] = "0.0182"
$arr[32,7] = 9
$arr[33,0] = 33
$arr[33,1] = "012605"
$arr[33,2] = "西藏东财中证证券保险领先指数A"
$arr[33,3] = "0.55"
$arr[33,4] = "94.99"
$arr[33,5] = "2.98"
$arr[33,6] = "0.0164"
$arr[33,7] = 9
$arr[34,0] = 34
$arr[34,1] = "013597"
$arr[34,2] = "招商中证全指证券公司指数（LOF）C"
$arr[34,3] = "0.39"
$arr[34,4] = "94.50"
$arr[34,5] = "3.30"
$arr[34,6] = "0.0129"
$arr[34,7] = 6
$arr[35,0] = 35
$arr[35,1] = "516200"
$arr[35,2] = "华安中证全指证券公司ETF"
$arr[35,3] = "0.32"
$arr[35,4] = "97.22"
$arr[35,5] = "3.08"
$arr[35,6] = "0.0099"
$arr[35,7] = 7
$arr[36,0] = 36
$arr[36,1] = "516980"
$arr[36,2] = "华富中证证券公司先锋策略ETF"
$arr[36,3] = "0.28"
$arr[36,4] = "99.02"
$arr[36,5] = "3.13"
$arr[36,6] = "0.0088"
$arr[36,7] = 6
$arr[37,0] = 37
$arr[37,1] = "015859"
$arr[37,2] = "宝盈国证证券龙头指数A"
$arr[37,3] = "0.12"
$arr[37,4] = "94.13"
$arr[37,5] = "6.85"
$arr[37,6] = "0.0082"
$arr[37,7] = 5
$arr[38,0] = 38
$arr[38,1] = "090011"
$arr[38,2] = "大成核心双动力混合"
$arr[38,3] = "0.24"
$arr[38,4] = "92.56"
$arr[38,5] = "3.34"
$arr[38,6] = "0.0080"
$arr[38,7] = 7
$arr[39,0] = 39
$arr[39,1] = "014984"
$arr[39,2] = "华安中证证券公司C"
$arr[39,3] = "0.18"
$arr[39,4] = "94.43"
$arr[39,5] = "3.29"
$arr[39,6] = "0.0059"
$arr[39,7] = 6
$arr[40,0] = 40
$arr[40,1] = "010404"
$arr[40,2] = "博道盛利6个月持有期混合"
$arr[40,3] = "1.10"
$arr[40,4] = "41.15"
$arr[40,5] = "0.46"
$arr[40,6] = "0.0051"
$arr[40,7] = 4
$arr[41,0] = 41
$arr[41,1] = "015860"
$arr[41,2] = "宝盈国证证券龙头指数C"
$arr[41,3] = "0.05"
$arr[41,4] = "94.13"
$arr[41,5] = "6.85"
$arr[41,6] = "0.0034"
$arr[41,7] = 5
$arr[42,0] = 42
$arr[42,1] = "015178"
$arr[42,2] = "申万菱信中证申万证券行业指数（LOF）C"
$arr[42,3] = "0.08"
$arr[42,4] = "93.19"
$arr[42,5] = "3.30"
$arr[42,6] = "0.0026"
$arr[42,7] = 6
$arr[43,0] = 43
$arr[43,1] = "015693"
$arr[43,2] = "鹏华中证800证券保险指数（LOF）C"
$arr[43,3] = "0.04"
$arr[43,4] = "94.43"
$arr[43,5] = "2.52"
$arr[43,6] = "0.0010"
$arr[43,7] = 9
$arr[44,0] = 44
$arr[44,1] = "015598"
$arr[44,2] = "国泰中证申万证券行业指数（LOF）C"
$arr[44,3] = "0.01"
$arr[44,4] = "93.42"
$arr[44,5] = "3.31"
$arr[44,6] = "0.0003"
$arr[44,7] = 6

$q3.Range("A2:H46").Value = $arr

# Re-apply the bold/border/center index-column style to the new A column
# (A2:A46), matching the same style used by A2 on the 总计 sheet and by the
# existing fund sheets.
$summary.Range("A2").Copy()
$q3.Range("A2:A46").PasteSpecial(-4122)

